$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item('LP1912')
$ws.Cells.Item(2,1).Value = 'Última actualización: 13:19:56'
$ws.Cells.Item(3,1).Value = 'Total filas: 265'
$ws.Cells.Item(121,1).Value = '08:21:50'
$ws.Cells.Item(121,2).Value = '09:41'
$ws.Cells.Item(121,3).Value = '215C_EL PATO'
$ws.Cells.Item(121,4).Value = 80
$ws.Cells.Item(121,5).Value = 'LP1912'
$ws.Cells.Item(122,1).Value = '09:38:09'
$ws.Cells.Item(122,2).Value = '09:41'
$ws.Cells.Item(122,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(122,4).Value = 3
$ws.Cells.Item(122,5).Value = 'LP1912'
$ws.Cells.Item(123,1).Value = '09:38:09'
$ws.Cells.Item(123,2).Value = '09:41'
$ws.Cells.Item(123,3).Value = '14_ABASTO'
$ws.Cells.Item(123,4).Value = 3
$ws.Cells.Item(123,5).Value = 'LP1912'
$ws.Cells.Item(204,1).Value = '12:01:50'
$ws.Cells.Item(204,2).Value = '12:34'
$ws.Cells.Item(204,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(204,4).Value = 33
$ws.Cells.Item(204,5).Value = 'LP1912'
$ws.Cells.Item(205,1).Value = '12:01:50'
$ws.Cells.Item(205,2).Value = '12:34'
$ws.Cells.Item(205,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(205,4).Value = 33
$ws.Cells.Item(205,5).Value = 'LP1912'
$ws.Cells.Item(229,1).Value = '13:19:56'
$ws.Cells.Item(229,2).Value = '13:20'
$ws.Cells.Item(229,3).Value = '26_HERNANDEZ'
$ws.Cells.Item(229,4).Value = 1
$ws.Cells.Item(229,5).Value = 'LP1912'
$ws.Cells.Item(230,1).Value = '11:48:04'
$ws.Cells.Item(230,2).Value = '13:20'
$ws.Cells.Item(230,3).Value = '10_OLMOS'
$ws.Cells.Item(230,4).Value = 92
$ws.Cells.Item(230,5).Value = 'LP1912'
$ws.Cells.Item(231,1).Value = '11:48:04'
$ws.Cells.Item(231,2).Value = '13:21'
$ws.Cells.Item(231,3).Value = '26_HERNANDEZ'
$ws.Cells.Item(231,4).Value = 93
$ws.Cells.Item(231,5).Value = 'LP1912'
$ws.Cells.Item(232,1).Value = '13:19:56'
$ws.Cells.Item(232,2).Value = '13:21'
$ws.Cells.Item(232,3).Value = '10_OLMOS'
$ws.Cells.Item(232,4).Value = 2
$ws.Cells.Item(232,5).Value = 'LP1912'
$ws.Cells.Item(233,1).Value = '13:19:56'
$ws.Cells.Item(233,2).Value = '13:22'
$ws.Cells.Item(233,3).Value = '14_ABASTO'
$ws.Cells.Item(233,4).Value = 3
$ws.Cells.Item(233,5).Value = 'LP1912'
$ws.Cells.Item(234,1).Value = '12:55:01'
$ws.Cells.Item(234,2).Value = '13:23'
$ws.Cells.Item(234,3).Value = '10_OLMOS'
$ws.Cells.Item(234,4).Value = 28
$ws.Cells.Item(234,5).Value = 'LP1912'
$ws.Cells.Item(235,1).Value = '13:19:56'
$ws.Cells.Item(235,2).Value = '13:26'
$ws.Cells.Item(235,3).Value = '15_ABASTO'
$ws.Cells.Item(235,4).Value = 7
$ws.Cells.Item(235,5).Value = 'LP1912'
$ws.Cells.Item(236,1).Value = '11:48:04'
$ws.Cells.Item(236,2).Value = '13:27'
$ws.Cells.Item(236,3).Value = '14_ABASTO'
$ws.Cells.Item(236,4).Value = 99
$ws.Cells.Item(236,5).Value = 'LP1912'
$ws.Cells.Item(237,1).Value = '13:19:56'
$ws.Cells.Item(237,2).Value = '13:31'
$ws.Cells.Item(237,3).Value = '10_OLMOS'
$ws.Cells.Item(237,4).Value = 12
$ws.Cells.Item(237,5).Value = 'LP1912'
$ws.Cells.Item(238,1).Value = '12:55:01'
$ws.Cells.Item(238,2).Value = '13:32'
$ws.Cells.Item(238,3).Value = '10_OLMOS'
$ws.Cells.Item(238,4).Value = 37
$ws.Cells.Item(238,5).Value = 'LP1912'
$ws.Cells.Item(239,1).Value = '13:19:56'
$ws.Cells.Item(239,2).Value = '13:34'
$ws.Cells.Item(239,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(239,4).Value = 15
$ws.Cells.Item(239,5).Value = 'LP1912'
$ws.Cells.Item(240,1).Value = '12:55:01'
$ws.Cells.Item(240,2).Value = '13:35'
$ws.Cells.Item(240,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(240,4).Value = 40
$ws.Cells.Item(240,5).Value = 'LP1912'
$ws.Cells.Item(241,1).Value = '11:48:04'
$ws.Cells.Item(241,2).Value = '13:36'
$ws.Cells.Item(241,3).Value = '15_ABASTO'
$ws.Cells.Item(241,4).Value = 108
$ws.Cells.Item(241,5).Value = 'LP1912'
$ws.Cells.Item(242,1).Value = '11:48:04'
$ws.Cells.Item(242,2).Value = '13:46'
$ws.Cells.Item(242,3).Value = '17_ROMERO'
$ws.Cells.Item(242,4).Value = 118
$ws.Cells.Item(242,5).Value = 'LP1912'
$ws.Cells.Item(243,1).Value = '13:19:56'
$ws.Cells.Item(243,2).Value = '13:46'
$ws.Cells.Item(243,3).Value = '16_SANTA ANA'
$ws.Cells.Item(243,4).Value = 27
$ws.Cells.Item(243,5).Value = 'LP1912'
$ws.Cells.Item(244,1).Value = '12:37:14'
$ws.Cells.Item(244,2).Value = '13:47'
$ws.Cells.Item(244,3).Value = '16_SANTA ANA'
$ws.Cells.Item(244,4).Value = 70
$ws.Cells.Item(244,5).Value = 'LP1912'
$ws.Cells.Item(245,1).Value = '13:19:56'
$ws.Cells.Item(245,2).Value = '13:50'
$ws.Cells.Item(245,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(245,4).Value = 31
$ws.Cells.Item(245,5).Value = 'LP1912'
$ws.Cells.Item(246,1).Value = '12:01:50'
$ws.Cells.Item(246,2).Value = '13:50'
$ws.Cells.Item(246,3).Value = '215A_EL PATO'
$ws.Cells.Item(246,4).Value = 109
$ws.Cells.Item(246,5).Value = 'LP1912'
$ws.Cells.Item(247,1).Value = '12:37:14'
$ws.Cells.Item(247,2).Value = '13:51'
$ws.Cells.Item(247,3).Value = '215A_EL PATO'
$ws.Cells.Item(247,4).Value = 74
$ws.Cells.Item(247,5).Value = 'LP1912'
$ws.Cells.Item(248,1).Value = '12:55:01'
$ws.Cells.Item(248,2).Value = '13:51'
$ws.Cells.Item(248,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(248,4).Value = 56
$ws.Cells.Item(248,5).Value = 'LP1912'
$ws.Cells.Item(249,1).Value = '12:01:50'
$ws.Cells.Item(249,2).Value = '13:52'
$ws.Cells.Item(249,3).Value = '10_OLMOS'
$ws.Cells.Item(249,4).Value = 111
$ws.Cells.Item(249,5).Value = 'LP1912'
$ws.Cells.Item(250,1).Value = '12:01:50'
$ws.Cells.Item(250,2).Value = '13:55'
$ws.Cells.Item(250,3).Value = '225_GOMEZ'
$ws.Cells.Item(250,4).Value = 114
$ws.Cells.Item(250,5).Value = 'LP1912'
$ws.Cells.Item(251,1).Value = '12:37:14'
$ws.Cells.Item(251,2).Value = '13:56'
$ws.Cells.Item(251,3).Value = '225_GOMEZ'
$ws.Cells.Item(251,4).Value = 79
$ws.Cells.Item(251,5).Value = 'LP1912'
$ws.Cells.Item(252,1).Value = '12:01:50'
$ws.Cells.Item(252,2).Value = '13:56'
$ws.Cells.Item(252,3).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(252,4).Value = 115
$ws.Cells.Item(252,5).Value = 'LP1912'
$ws.Cells.Item(253,1).Value = '12:37:14'
$ws.Cells.Item(253,2).Value = '13:57'
$ws.Cells.Item(253,3).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(253,4).Value = 80
$ws.Cells.Item(253,5).Value = 'LP1912'
$ws.Cells.Item(254,1).Value = '13:19:56'
$ws.Cells.Item(254,2).Value = '14:04'
$ws.Cells.Item(254,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(254,4).Value = 45
$ws.Cells.Item(254,5).Value = 'LP1912'
$ws.Cells.Item(255,1).Value = '12:37:14'
$ws.Cells.Item(255,2).Value = '14:04'
$ws.Cells.Item(255,3).Value = '17_ROMERO'
$ws.Cells.Item(255,4).Value = 87
$ws.Cells.Item(255,5).Value = 'LP1912'
$ws.Cells.Item(256,1).Value = '13:19:56'
$ws.Cells.Item(256,2).Value = '14:06'
$ws.Cells.Item(256,3).Value = '16_SANTA ANA'
$ws.Cells.Item(256,4).Value = 47
$ws.Cells.Item(256,5).Value = 'LP1912'
$ws.Cells.Item(257,1).Value = '12:55:01'
$ws.Cells.Item(257,2).Value = '14:07'
$ws.Cells.Item(257,3).Value = '16_SANTA ANA'
$ws.Cells.Item(257,4).Value = 72
$ws.Cells.Item(257,5).Value = 'LP1912'
$ws.Cells.Item(258,1).Value = '13:19:56'
$ws.Cells.Item(258,2).Value = '14:16'
$ws.Cells.Item(258,3).Value = '27_EL RETIRO'
$ws.Cells.Item(258,4).Value = 57
$ws.Cells.Item(258,5).Value = 'LP1912'
$ws.Cells.Item(259,1).Value = '12:37:14'
$ws.Cells.Item(259,2).Value = '14:17'
$ws.Cells.Item(259,3).Value = '27_EL RETIRO'
$ws.Cells.Item(259,4).Value = 100
$ws.Cells.Item(259,5).Value = 'LP1912'
$ws.Cells.Item(260,1).Value = '13:19:56'
$ws.Cells.Item(260,2).Value = '14:19'
$ws.Cells.Item(260,3).Value = '215C_EL PATO'
$ws.Cells.Item(260,4).Value = 60
$ws.Cells.Item(260,5).Value = 'LP1912'
$ws.Cells.Item(261,1).Value = '12:37:14'
$ws.Cells.Item(261,2).Value = '14:20'
$ws.Cells.Item(261,3).Value = '215C_EL PATO'
$ws.Cells.Item(261,4).Value = 103
$ws.Cells.Item(261,5).Value = 'LP1912'
$ws.Cells.Item(262,1).Value = '13:19:56'
$ws.Cells.Item(262,2).Value = '14:20'
$ws.Cells.Item(262,3).Value = '26_HERNANDEZ'
$ws.Cells.Item(262,4).Value = 61
$ws.Cells.Item(262,5).Value = 'LP1912'
$ws.Cells.Item(263,1).Value = '12:37:14'
$ws.Cells.Item(263,2).Value = '14:21'
$ws.Cells.Item(263,3).Value = '26_HERNANDEZ'
$ws.Cells.Item(263,4).Value = 104
$ws.Cells.Item(263,5).Value = 'LP1912'
$ws.Cells.Item(264,1).Value = '13:19:56'
$ws.Cells.Item(264,2).Value = '14:49'
$ws.Cells.Item(264,3).Value = '14_ABASTO'
$ws.Cells.Item(264,4).Value = 90
$ws.Cells.Item(264,5).Value = 'LP1912'
$ws.Cells.Item(265,1).Value = '12:55:01'
$ws.Cells.Item(265,2).Value = '14:50'
$ws.Cells.Item(265,3).Value = '14_ABASTO'
$ws.Cells.Item(265,4).Value = 115
$ws.Cells.Item(265,5).Value = 'LP1912'
$ws.Cells.Item(266,1).Value = '13:19:56'
$ws.Cells.Item(266,2).Value = '14:56'
$ws.Cells.Item(266,3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(266,4).Value = 97
$ws.Cells.Item(266,5).Value = 'LP1912'
$ws.Cells.Item(267,1).Value = '13:19:56'
$ws.Cells.Item(267,2).Value = '14:58'
$ws.Cells.Item(267,3).Value = '215B_EL PATO'
$ws.Cells.Item(267,4).Value = 99
$ws.Cells.Item(267,5).Value = 'LP1912'
$ws.Cells.Item(268,1).Value = '13:19:56'
$ws.Cells.Item(268,2).Value = '15:00'
$ws.Cells.Item(268,3).Value = '81_EL PELIGRO'
$ws.Cells.Item(268,4).Value = 101
$ws.Cells.Item(268,5).Value = 'LP1912'
$ws.Cells.Item(269,1).Value = '13:19:56'
$ws.Cells.Item(269,2).Value = '15:04'
$ws.Cells.Item(269,3).Value = '10_OLMOS'
$ws.Cells.Item(269,4).Value = 105
$ws.Cells.Item(269,5).Value = 'LP1912'
$ws.Cells.Item(270,1).Value = '13:19:56'
$ws.Cells.Item(270,2).Value = '15:13'
$ws.Cells.Item(270,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(270,4).Value = 114
$ws.Cells.Item(270,5).Value = 'LP1912'

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item('LP1912-215')
$ws.Cells.Item(2,1).Value = 'Última actualización: 13:19:56'
$ws.Cells.Item(3,1).Value = 'Total filas: 31'
$ws.Cells.Item(34,1).Value = '13:19:56'
$ws.Cells.Item(34,2).Value = '14:19'
$ws.Cells.Item(34,3).Value = '215C_EL PATO'
$ws.Cells.Item(34,4).Value = 60
$ws.Cells.Item(34,5).Value = 'LP1912'
$ws.Cells.Item(35,1).Value = '12:37:14'
$ws.Cells.Item(35,2).Value = '14:20'
$ws.Cells.Item(35,3).Value = '215C_EL PATO'
$ws.Cells.Item(35,4).Value = 103
$ws.Cells.Item(35,5).Value = 'LP1912'
$ws.Cells.Item(36,1).Value = '13:19:56'
$ws.Cells.Item(36,2).Value = '14:58'
$ws.Cells.Item(36,3).Value = '215B_EL PATO'
$ws.Cells.Item(36,4).Value = 99
$ws.Cells.Item(36,5).Value = 'LP1912'

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item('6203-6173')
$ws.Cells.Item(2,1).Value = 'Última actualización: 13:19:56'
$ws.Cells.Item(3,1).Value = 'Total filas: 42'
$ws.Cells.Item(44,1).Value = '13:19:56'
$ws.Cells.Item(44,2).Value = '14:08'
$ws.Cells.Item(44,3).Value = '215A_LA PLATA'
$ws.Cells.Item(44,4).Value = 49
$ws.Cells.Item(44,5).Value = 'L6173'
$ws.Cells.Item(45,1).Value = '12:37:14'
$ws.Cells.Item(45,2).Value = '14:09'
$ws.Cells.Item(45,3).Value = '215A_LA PLATA'
$ws.Cells.Item(45,4).Value = 92
$ws.Cells.Item(45,5).Value = 'L6173'
$ws.Cells.Item(46,1).Value = '13:19:56'
$ws.Cells.Item(46,2).Value = '14:52'
$ws.Cells.Item(46,3).Value = '215D_LA PLATA'
$ws.Cells.Item(46,4).Value = 93
$ws.Cells.Item(46,5).Value = 'L6203'
$ws.Cells.Item(47,1).Value = '12:55:01'
$ws.Cells.Item(47,2).Value = '14:53'
$ws.Cells.Item(47,3).Value = '215D_LA PLATA'
$ws.Cells.Item(47,4).Value = 118
$ws.Cells.Item(47,5).Value = 'L6203'
